$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Duplicate the formatting of the last existing data row (row 43) onto the
# two new rows so the new cells keep the same style indices (date format,
# currency-like number format, integer format) instead of Excel creating
# brand-new style entries.
$ws.Range("A43:F43").Copy()
$ws.Range("A44:F44").PasteSpecial(-4122)
$ws.Range("A45:F45").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 44: 四方坪站
$ws.Cells.Item(44, 1).Value = 46044
$ws.Cells.Item(44, 2).Value = "四方坪站"
$ws.Cells.Item(44, 3).Value = 15956.26
$ws.Cells.Item(44, 4).Value = 12826.05
$ws.Cells.Item(44, 5).Value = 5217.01
$ws.Cells.Item(44, 6).Value = 638

# Row 45: 高岭站
$ws.Cells.Item(45, 1).Value = 46044
$ws.Cells.Item(45, 2).Value = "高岭站"
$ws.Cells.Item(45, 3).Value = 5739.92
$ws.Cells.Item(45, 4).Value = 4920.21
$ws.Cells.Item(45, 5).Value = 1523.49
$ws.Cells.Item(45, 6).Value = 209

# Mirror the active-cell selection shift recorded in the diff (I42 -> I43)
$ws.Range("I43").Select()
